# Weekly update: a new daily price record is inserted at row 301 ("Fruta /
# hortaliza, semanal"). Excel shifts all the existing rows (301-317) down by
# one (they become 302-318), carrying their values and formatting with them,
# and the sheet's used range grows from A1:T317 to A1:T318.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 301; everything below moves
# down one row automatically (values, number formats, etc.).
$ws.Rows.Item(301).Insert()

# Populate the newly inserted row 301 with the new weekly record.
$ws.Range("A301").Value = 10
$ws.Range("B301").Value = 'Vega Modelo de Temuco'
$ws.Range("C301").Value = 'La Araucanía'
$ws.Range("D301").Value = 44516
$ws.Range("E301").Value = 9
$ws.Range("F301").Value = 'Fruta'
$ws.Range("G301").Value = 100108
$ws.Range("H301").Value = 'Tropicales y subtropicales'
$ws.Range("I301").Value = 100108005
$ws.Range("J301").Value = 'Piña'
$ws.Range("K301").Value = 'Caramelo'
$ws.Range("L301").Value = 'Primera'
$ws.Range("M301").Value = 95
$ws.Range("N301").Value = 19000
$ws.Range("O301").Value = 20000
$ws.Range("P301").Value = 19474
$ws.Range("Q301").Value = '$/caja 12 unidades'
$ws.Range("R301").Value = 'Ecuador'
$ws.Range("S301").Value = 1623
$ws.Range("T301").Value = 12
